# AddTurmaEAlunos.xlsx - "Alteração do type date nos ficheiros Excel"
#
# Turma sheet:   Edição (B2) "99.99" -> "12.20" (kept as text)
# Formandos sheet: Nº de Formando (A2) "T0000000" -> "T0000001"
#                   Email (C2, hyperlink display text) updated to match
#                   Data de Nascimento (D2) gets an explicit dd/mm/yyyy date format
# Active sheet switches from "Turma" to "Formandos", selections updated too.

$wb = $excel.ActiveWorkbook

$wsTurma = $wb.Worksheets.Item("Turma")
$wsFormandos = $wb.Worksheets.Item("Formandos")

# --- Turma sheet: Edição value update (stored as text, like "99.99" was) ---
$wsTurma.Range("B2").NumberFormat = "@"
$wsTurma.Range("B2").Value = "12.20"

# --- Formandos sheet: student number + derived email ---
$wsFormandos.Range("A2").Value = "T0000001"
$wsFormandos.Range("C2").Value = "Nome.Apelido.T0001000@atec.pt"

# --- Formandos sheet: birth date column now uses an explicit dd/mm/yyyy format ---
$wsFormandos.Range("D2").NumberFormat = "dd/mm/yyyy;@"

# --- View state: Formandos becomes the active/selected tab ---
$wsTurma.Activate() | Out-Null
$wsTurma.Range("C2").Select() | Out-Null
$wsFormandos.Activate() | Out-Null
$wsFormandos.Range("A6").Select() | Out-Null
